$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 rows at row 17 (pushes old row17 data to row20, footer rows 22/23 to 25/26)
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()

# 2. Copy the formatting of row 16 onto the freshly inserted rows 17-19
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# 3. Clear the old shared-string references so the shared-strings table can be rebuilt
#    in the exact order the new content requires.
$ws.Range("C16:E16").ClearContents()
$ws.Range("C20:E20").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("B13").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("H15:J15").ClearContents()
$ws.Range("B25:C26").ClearContents()
$ws.Range("H25:J26").ClearContents()

# 4. Re-populate everything, in the precise order needed to reproduce the target
#    shared-strings table ordering.

# Row 16 - MACYORI OYOLA VEGA / periodo 2210
$ws.Cells.Item(16, 3).Value = "1047382783"
$ws.Cells.Item(16, 4).Value = "MACYORI OYOLA VEGA"
$ws.Cells.Item(16, 5).Value = "2210"
$ws.Cells.Item(16, 6).Value = 58506
$ws.Cells.Item(16, 7).Value = 1462644

# Row 17 - MACYORI OYOLA VEGA / periodo 2209
$ws.Cells.Item(17, 2).Value = "CC"
$ws.Cells.Item(17, 3).Value = "1047382783"
$ws.Cells.Item(17, 4).Value = "MACYORI OYOLA VEGA"
$ws.Cells.Item(17, 5).Value = "2209"
$ws.Cells.Item(17, 6).Value = 58506
$ws.Cells.Item(17, 7).Value = 1462644

# Row 18 - MILFRE LILIANA MOSQUERA IBARGUEN / periodo 2506
$ws.Cells.Item(18, 2).Value = "CC"
$ws.Cells.Item(18, 3).Value = "35851218"
$ws.Cells.Item(18, 4).Value = "MILFRE LILIANA MOSQUERA IBARGUEN"
$ws.Cells.Item(18, 5).Value = "2506"
$ws.Cells.Item(18, 6).Value = 45552
$ws.Cells.Item(18, 7).Value = 1423500

# Row 19 - MAIRA ALEJANDRA JURADO FERNANDEZ / periodo 2302
$ws.Cells.Item(19, 2).Value = "CC"
$ws.Cells.Item(19, 3).Value = "1052958948"
$ws.Cells.Item(19, 4).Value = "MAIRA ALEJANDRA JURADO FERNANDEZ"
$ws.Cells.Item(19, 5).Value = "2302"
$ws.Cells.Item(19, 6).Value = 43002
$ws.Cells.Item(19, 7).Value = 1160000

# Row 20 - LILIEN TAINA BARRIOS OLIVO / periodo 2105
$ws.Cells.Item(20, 3).Value = "1049929498"
$ws.Cells.Item(20, 4).Value = "LILIEN TAINA BARRIOS OLIVO"
$ws.Cells.Item(20, 5).Value = "2105"

# Header/labels
$ws.Range("D2").Value = "ESTADO DE CUENTA"
$ws.Range("B7").Value = "RAZON SOCIAL:"
$ws.Range("B11").Value = "VALOR MORA"
$ws.Cells.Item(11, 5).Value = 214872
$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Cells.Item(13, 3).Value = 4
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Cells.Item(13, 6).Value = 5

# Table header row 15
$ws.Cells.Item(15, 9).Value = "Novedad de Retiro"
$ws.Cells.Item(15, 8).Value = "Novedad de Ingreso"
$ws.Cells.Item(15, 10).Value = "Observaciones"

# Footer rows 25/26
$ws.Cells.Item(26, 2).Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Cells.Item(25, 2).Value = "___________________________________"
$ws.Cells.Item(26, 8).Value = "FIRMA DEL REPRESENTANTE LEGAL"
$ws.Cells.Item(25, 8).Value = "___________________________________"
